$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.703.51"
$ws.Range("E2").Value = "  -3.36%  "
$ws.Range("D3").Value = "3.816.14"
$ws.Range("E3").Value = "  -2.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D7").Value = "3.804.79"
$ws.Range("E7").Value = "  -3.26%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.49%  "
$ws.Range("D15").Value = "4.473.42"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").Value = "3.837.66"
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("D17").Value = "67.860.25"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.81%  "
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "463.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000160"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").Value = "3.980.58"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("D36").Value = "3.794.43"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("E37").Value = "  -4.24%  "
$ws.Range("E38").Value = "  +9.88%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.139"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.312"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.73%  "
$ws.Range("B45").Value = "FLOKI"
$ws.Range("C45").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000296"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.13%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "418.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.86%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.15%  "
